$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet view changes
$ws.Range("H14").Select()
$ws.Application.ActiveWindow.ScrollRow = 25

Write-Host "test"
